$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (far outside the table's used range) that evaluates to the
# literal text "1". Using a formula keeps the result as a genuine text
# value (Excel's normal "smart" numeric auto-conversion only kicks in for
# literal Range.Value assignments, not for cached formula results), and
# copying it elsewhere via PasteSpecial(xlPasteValues) overwrites just the
# destination's value while preserving the destination's existing style.
$ws.Range("Z20").Formula = "=""1"""

# --- Row 5: existing WAN device entry gets a new IP / hostname / CPU-UTILS value ---
$ws.Range("B5").Value = "192.168.0.106"
$ws.Range("C5").Value = "lab-3725"
$ws.Range("Z20").Copy()
$ws.Range("G5").PasteSpecial(-4163)

# --- Rows 6-9: four new WAN device entries, matching row 5's look & feel ---
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G9").PasteSpecial(-4122)

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "192.168.0.107"
$ws.Range("C6").Value = "lab-3660"
$ws.Range("D6").Value = "YES"
$ws.Range("E6").Value = "YES"
$ws.Range("F6").Value = "YES"
$ws.Range("Z20").Copy()
$ws.Range("G6").PasteSpecial(-4163)

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "192.168.0.108"
$ws.Range("C7").Value = "lab-3640"
$ws.Range("D7").Value = "YES"
$ws.Range("E7").Value = "YES"
$ws.Range("F7").Value = "YES"
$ws.Range("Z20").Copy()
$ws.Range("G7").PasteSpecial(-4163)

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "192.168.0.109"
$ws.Range("C8").Value = "lab-3745"
$ws.Range("D8").Value = "YES"
$ws.Range("E8").Value = "YES"
$ws.Range("F8").Value = "YES"
$ws.Range("Z20").Copy()
$ws.Range("G8").PasteSpecial(-4163)

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "192.168.0.110"
$ws.Range("C9").Value = "lab-2691"
$ws.Range("D9").Value = "YES"
$ws.Range("E9").Value = "YES"
$ws.Range("F9").Value = "YES"
$ws.Range("Z20").Copy()
$ws.Range("G9").PasteSpecial(-4163)

# --- Clean up the helper cell so it doesn't linger in the saved sheet ---
$ws.Range("Z20").Clear()

Write-Host "NE-AUDIT: rows 5-9 updated (4 new WAN devices added)"
